$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells (row 1) - values
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold, border, centered) from H1 to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data rows 2-7 for columns I and J
$values = @{
    2 = 9
    3 = 8
    4 = 5
    5 = 7
    6 = 3
    7 = 4
}

foreach ($row in $values.Keys) {
    $val = $values[$row]
    $ws.Cells.Item($row, 9).Value = $val   # Column I
    $ws.Cells.Item($row, 10).Value = $val  # Column J
}
